$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C ("Förändrad") for rows 2 through 185 from serial date 45172 to 45175
$ws.Range("C2:C185").Value = 45175
